$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 is a brand-new annotation row. Column B there holds the text "3"
# (same look as the old B37 had before it got normalized to a real number).
# Copy the still-intact old B37 cell down to B38 first, so the text type is
# preserved faithfully without Excel tagging it with a "number stored as
# text" quote-prefix style.
$ws.Range("B37").Copy()
$ws.Range("B38").PasteSpecial()

# Now correct B37 itself: store it as a genuine number (matching every
# other row in the politeness_score column).
$ws.Range("B37").Value = 3

# Fill in the rest of the new row 38.
$ws.Range("A38").Value = "Ruilin"
$ws.Range("C38").Value = "无"
$ws.Range("D38").Value = "APC"
$ws.Range("E38").Value = "MET"
$ws.Range("F38").Value = "d0296b92-10f5-497e-8726-aae675ac805b"
$ws.Range("G38").Value = "rJl3yM-Ab_annotated.xlsx"
$ws.Range("H38").Value = "The new method is motivated well and departs from prior work."
